$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows before row 1154, shifting existing rows 1154-1206 down to 1160-1212
$ws.Rows("1154:1159").Insert()

# Row 1154
$ws.Cells.Item(1154, 1).Value = 9
$ws.Cells.Item(1154, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(1154, 3).Value = 'Metropolitana'
$ws.Cells.Item(1154, 4).Value = 44509
$ws.Cells.Item(1154, 5).Value = 13
$ws.Cells.Item(1154, 6).Value = 'Fruta'
$ws.Cells.Item(1154, 7).Value = 100102
$ws.Cells.Item(1154, 8).Value = 'Cítricos'
$ws.Cells.Item(1154, 9).Value = 100102003
$ws.Cells.Item(1154, 10).Value = 'Limón'
$ws.Cells.Item(1154, 11).Value = 'Sin especificar'
$ws.Cells.Item(1154, 12).Value = '1a amarillo'
$ws.Cells.Item(1154, 13).Value = 400
$ws.Cells.Item(1154, 14).Value = 10000
$ws.Cells.Item(1154, 15).Value = 10000
$ws.Cells.Item(1154, 16).Value = 10000
$ws.Cells.Item(1154, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(1154, 18).Value = 'Provincia del Elquí'
$ws.Cells.Item(1154, 19).Value = 556
$ws.Cells.Item(1154, 20).Value = 18

# Row 1155
$ws.Cells.Item(1155, 1).Value = 9
$ws.Cells.Item(1155, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(1155, 3).Value = 'Metropolitana'
$ws.Cells.Item(1155, 4).Value = 44509
$ws.Cells.Item(1155, 5).Value = 13
$ws.Cells.Item(1155, 6).Value = 'Fruta'
$ws.Cells.Item(1155, 7).Value = 100102
$ws.Cells.Item(1155, 8).Value = 'Cítricos'
$ws.Cells.Item(1155, 9).Value = 100102003
$ws.Cells.Item(1155, 10).Value = 'Limón'
$ws.Cells.Item(1155, 11).Value = 'Sin especificar'
$ws.Cells.Item(1155, 12).Value = '1a amarillo'
$ws.Cells.Item(1155, 13).Value = 680
$ws.Cells.Item(1155, 14).Value = 6000
$ws.Cells.Item(1155, 15).Value = 6500
$ws.Cells.Item(1155, 16).Value = 6257
$ws.Cells.Item(1155, 17).Value = '$/malla 18 kilos'
$ws.Cells.Item(1155, 18).Value = 'Provincia del Elquí'
$ws.Cells.Item(1155, 19).Value = 348
$ws.Cells.Item(1155, 20).Value = 18

# Row 1156
$ws.Cells.Item(1156, 1).Value = 9
$ws.Cells.Item(1156, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(1156, 3).Value = 'Metropolitana'
$ws.Cells.Item(1156, 4).Value = 44509
$ws.Cells.Item(1156, 5).Value = 13
$ws.Cells.Item(1156, 6).Value = 'Fruta'
$ws.Cells.Item(1156, 7).Value = 100102
$ws.Cells.Item(1156, 8).Value = 'Cítricos'
$ws.Cells.Item(1156, 9).Value = 100102003
$ws.Cells.Item(1156, 10).Value = 'Limón'
$ws.Cells.Item(1156, 11).Value = 'Sin especificar'
$ws.Cells.Item(1156, 12).Value = '1a amarillo'
$ws.Cells.Item(1156, 13).Value = 330
$ws.Cells.Item(1156, 14).Value = 6000
$ws.Cells.Item(1156, 15).Value = 6000
$ws.Cells.Item(1156, 16).Value = 6000
$ws.Cells.Item(1156, 17).Value = '$/malla 18 kilos'
$ws.Cells.Item(1156, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(1156, 19).Value = 333
$ws.Cells.Item(1156, 20).Value = 18

# Row 1157
$ws.Cells.Item(1157, 1).Value = 9
$ws.Cells.Item(1157, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(1157, 3).Value = 'Metropolitana'
$ws.Cells.Item(1157, 4).Value = 44509
$ws.Cells.Item(1157, 5).Value = 13
$ws.Cells.Item(1157, 6).Value = 'Fruta'
$ws.Cells.Item(1157, 7).Value = 100102
$ws.Cells.Item(1157, 8).Value = 'Cítricos'
$ws.Cells.Item(1157, 9).Value = 100102003
$ws.Cells.Item(1157, 10).Value = 'Limón'
$ws.Cells.Item(1157, 11).Value = 'Sin especificar'
$ws.Cells.Item(1157, 12).Value = '2a amarillo'
$ws.Cells.Item(1157, 13).Value = 500
$ws.Cells.Item(1157, 14).Value = 5000
$ws.Cells.Item(1157, 15).Value = 5500
$ws.Cells.Item(1157, 16).Value = 5280
$ws.Cells.Item(1157, 17).Value = '$/malla 18 kilos'
$ws.Cells.Item(1157, 18).Value = 'Provincia del Elquí'
$ws.Cells.Item(1157, 19).Value = 293
$ws.Cells.Item(1157, 20).Value = 18

# Row 1158
$ws.Cells.Item(1158, 1).Value = 9
$ws.Cells.Item(1158, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(1158, 3).Value = 'Metropolitana'
$ws.Cells.Item(1158, 4).Value = 44509
$ws.Cells.Item(1158, 5).Value = 13
$ws.Cells.Item(1158, 6).Value = 'Fruta'
$ws.Cells.Item(1158, 7).Value = 100102
$ws.Cells.Item(1158, 8).Value = 'Cítricos'
$ws.Cells.Item(1158, 9).Value = 100102003
$ws.Cells.Item(1158, 10).Value = 'Limón'
$ws.Cells.Item(1158, 11).Value = 'Sin especificar'
$ws.Cells.Item(1158, 12).Value = '2a amarillo'
$ws.Cells.Item(1158, 13).Value = 480
$ws.Cells.Item(1158, 14).Value = 5000
$ws.Cells.Item(1158, 15).Value = 5000
$ws.Cells.Item(1158, 16).Value = 5000
$ws.Cells.Item(1158, 17).Value = '$/malla 18 kilos'
$ws.Cells.Item(1158, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(1158, 19).Value = 278
$ws.Cells.Item(1158, 20).Value = 18

# Row 1159
$ws.Cells.Item(1159, 1).Value = 9
$ws.Cells.Item(1159, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(1159, 3).Value = 'Metropolitana'
$ws.Cells.Item(1159, 4).Value = 44509
$ws.Cells.Item(1159, 5).Value = 13
$ws.Cells.Item(1159, 6).Value = 'Fruta'
$ws.Cells.Item(1159, 7).Value = 100102
$ws.Cells.Item(1159, 8).Value = 'Cítricos'
$ws.Cells.Item(1159, 9).Value = 100102003
$ws.Cells.Item(1159, 10).Value = 'Limón'
$ws.Cells.Item(1159, 11).Value = 'Sin especificar'
$ws.Cells.Item(1159, 12).Value = '3a amarillo'
$ws.Cells.Item(1159, 13).Value = 420
$ws.Cells.Item(1159, 14).Value = 3500
$ws.Cells.Item(1159, 15).Value = 3500
$ws.Cells.Item(1159, 16).Value = 3500
$ws.Cells.Item(1159, 17).Value = '$/malla 18 kilos'
$ws.Cells.Item(1159, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(1159, 19).Value = 194
$ws.Cells.Item(1159, 20).Value = 18
